# iOS开发进度.xlsx - "finished user account setting view"
#
# Updates the progress tracker:
#  - Several already-complete rows (服务详请/门店详情/公益列表/活动列表/
#    门店-攻略-商品-资讯评论 and 短信) have their "预计完成时间" reset to "--"
#    now that they're fully done (same visual style as the other "--" rows).
#  - 个人中心 (user account / personal center) is now 100% complete; its
#    note about missing design assets is cleared.
#  - 地图 / 支付 / 登录注册 / 消息 / 引导页 get refreshed target dates / notes.
#  - A brand-new row is appended for the 分享 (Share) module.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# --- Rows that are now fully finished: predicted date column becomes "--" ---
# (reuse the style already used by the other "--" cells, e.g. C3:C5)
"C6", "C7", "C8", "C9", "C10", "C15" | ForEach-Object {
    Copy-CellFormat "C3" $_
    $ws.Range($_).Value2 = "--"
}
$excel.CutCopyMode = $false

# --- Row 11: 个人中心 (Personal center / account settings) -> finished ---
Copy-CellFormat "B3" "B11"
$ws.Range("B11").Value2 = 1
$excel.CutCopyMode = $false
$ws.Range("C11").Value2 = "11.13"
$ws.Range("D11").ClearContents() | Out-Null

# --- Row 12: 地图 (Map) -> new target date ---
$ws.Range("C12").Value2 = "11.17"

# --- Row 13: 支付 (Payment) -> target date refreshed ---
$ws.Range("C13").Value2 = "11.15"

# --- Row 14: 登录/注册 (Login/Register) -> progress + date updated ---
$ws.Range("B14").Value2 = 0.5
$ws.Range("C14").Value2 = "11.14"
$ws.Range("D14").Value2 = "待接通接口"

# --- Row 16: 消息 (Message) -> target date refreshed ---
$ws.Range("C16").Value2 = "11.16"

# --- Row 17: 引导页 (Onboarding) -> date refreshed + note updated ---
$ws.Range("C17").Value2 = "11.15"
$ws.Range("D17").Value2 = "待切图"

# --- New row 18: 分享 (Share) module ---
$ws.Range("A18").Value2 = "分享"
$ws.Range("B18").NumberFormat = "0.00%"
$ws.Range("B18").Value2 = 0.2
$ws.Range("C18").Value2 = "11.18"

$ws.Range("D18").Select() | Out-Null
